# Adds rows to "Programacao" and "Planilha" sheets, and updates several
# cell values on "Descarga do Sal", per the target diff.
$wb = $excel.ActiveWorkbook

# --- Sheet "Programacao": add rows 43-46 ---
$ws1 = $wb.Worksheets.Item("Programacao")
$ws1.Cells.Item(43,1).Value = "das"
$ws1.Cells.Item(43,2).Value = "das"
$ws1.Cells.Item(43,3).Value = "das"
$ws1.Cells.Item(43,4).Value = "das"
$ws1.Cells.Item(43,5).Value = "das"
$ws1.Cells.Item(43,6).Value = "NORSAL"
$ws1.Cells.Item(43,7).Value = 516
$ws1.Cells.Item(43,8).Value = "das"
$ws1.Cells.Item(43,9).Value = "das"
$ws1.Cells.Item(43,10).Value = "SAL REFINADO Selecione uma opção"
$ws1.Cells.Item(43,11).Value = "das"

$ws1.Cells.Item(44,1).Value = "das"
$ws1.Cells.Item(44,2).Value = "das"
$ws1.Cells.Item(44,3).Value = "dasd"
$ws1.Cells.Item(44,4).Value = "das"
$ws1.Cells.Item(44,5).Value = "das"
$ws1.Cells.Item(44,6).Value = "NORSAL"
$ws1.Cells.Item(44,7).Value = 265
$ws1.Cells.Item(44,8).Value = "das"
$ws1.Cells.Item(44,9).Value = "das"
$ws1.Cells.Item(44,10).Value = "SAL REFINADO Selecione uma opção"
$ws1.Cells.Item(44,11).Value = "das"

$ws1.Cells.Item(45,1).Value = "das"
$ws1.Cells.Item(45,2).Value = "das"
$ws1.Cells.Item(45,3).Value = "dasd"
$ws1.Cells.Item(45,4).Value = "das"
$ws1.Cells.Item(45,5).Value = "das"
$ws1.Cells.Item(45,6).Value = "NORSAL"
$ws1.Cells.Item(45,7).Value = 265
$ws1.Cells.Item(45,8).Value = "das"
$ws1.Cells.Item(45,9).Value = "das"
$ws1.Cells.Item(45,10).Value = "SAL REFINADO Selecione uma opção"
$ws1.Cells.Item(45,11).Value = "das"

$ws1.Cells.Item(46,1).Value = "das"
$ws1.Cells.Item(46,2).Value = "das"
$ws1.Cells.Item(46,3).Value = "dasd"
$ws1.Cells.Item(46,4).Value = "das"
$ws1.Cells.Item(46,5).Value = "das"
$ws1.Cells.Item(46,6).Value = "NORSAL"
$ws1.Cells.Item(46,7).Value = 265
$ws1.Cells.Item(46,8).Value = "das"
$ws1.Cells.Item(46,9).Value = "das"
$ws1.Cells.Item(46,10).Value = "SAL REFINADO Selecione uma opção"
$ws1.Cells.Item(46,11).Value = "das"

# --- Sheet "Planilha": add rows 68-74 ---
$ws2 = $wb.Worksheets.Item("Planilha")
$ws2.Cells.Item(68,1).Value = "ENTRADA"
$ws2.Cells.Item(68,2).Value = "das"
$ws2.Cells.Item(68,3).Value = "das"
$ws2.Cells.Item(68,4).Value = "das"
$ws2.Cells.Item(68,5).Value = "das"
$ws2.Cells.Item(68,6).Value = "SAL REFINADO"
$ws2.Cells.Item(68,7).Value = "Selecione uma opção"
$ws2.Cells.Item(68,8).Value = "NORSAL"
$ws2.Cells.Item(68,9).Value = "das"
$ws2.Cells.Item(68,10).Value = "das"
$ws2.Cells.Item(68,11).Value = 561
$ws2.Cells.Item(68,12).Value = "das"
$ws2.Cells.Item(68,13).Value = "das"
$ws2.Cells.Item(68,14).Value = 516

$ws2.Cells.Item(69,1).Value = "ENTRADA"
$ws2.Cells.Item(69,2).Value = "das"
$ws2.Cells.Item(69,3).Value = "das"
$ws2.Cells.Item(69,4).Value = "das"
$ws2.Cells.Item(69,5).Value = "das"
$ws2.Cells.Item(69,6).Value = "SAL REFINADO"
$ws2.Cells.Item(69,7).Value = "Selecione uma opção"
$ws2.Cells.Item(69,8).Value = "NORSAL"
$ws2.Cells.Item(69,9).Value = "das"
$ws2.Cells.Item(69,10).Value = "das"
$ws2.Cells.Item(69,11).Value = 68
$ws2.Cells.Item(69,12).Value = "das"
$ws2.Cells.Item(69,13).Value = "das"
$ws2.Cells.Item(69,14).Value = 265

$ws2.Cells.Item(70,1).Value = "ENTRADA"
$ws2.Cells.Item(70,2).Value = "das"
$ws2.Cells.Item(70,3).Value = "das"
$ws2.Cells.Item(70,4).Value = "das"
$ws2.Cells.Item(70,5).Value = "das"
$ws2.Cells.Item(70,6).Value = "SAL REFINADO"
$ws2.Cells.Item(70,7).Value = "Selecione uma opção"
$ws2.Cells.Item(70,8).Value = "NORSAL"
$ws2.Cells.Item(70,9).Value = "das"
$ws2.Cells.Item(70,10).Value = "das"
$ws2.Cells.Item(70,11).Value = 68
$ws2.Cells.Item(70,12).Value = "das"
$ws2.Cells.Item(70,13).Value = "das"
$ws2.Cells.Item(70,14).Value = 265

$ws2.Cells.Item(71,1).Value = "ENTRADA"
$ws2.Cells.Item(71,2).Value = "das"
$ws2.Cells.Item(71,3).Value = "das"
$ws2.Cells.Item(71,4).Value = "das"
$ws2.Cells.Item(71,5).Value = "das"
$ws2.Cells.Item(71,6).Value = "SAL REFINADO"
$ws2.Cells.Item(71,7).Value = "Selecione uma opção"
$ws2.Cells.Item(71,8).Value = "NORSAL"
$ws2.Cells.Item(71,9).Value = "das"
$ws2.Cells.Item(71,10).Value = "das"
$ws2.Cells.Item(71,11).Value = 165
$ws2.Cells.Item(71,12).Value = "das"
$ws2.Cells.Item(71,13).Value = "das"
$ws2.Cells.Item(71,14).Value = 65

$ws2.Cells.Item(72,1).Value = "ENTRADA"
$ws2.Cells.Item(72,2).Value = "das"
$ws2.Cells.Item(72,3).Value = "das"
$ws2.Cells.Item(72,4).Value = "das"
$ws2.Cells.Item(72,5).Value = "das"
$ws2.Cells.Item(72,6).Value = "SAL REFINADO"
$ws2.Cells.Item(72,7).Value = "Selecione uma opção"
$ws2.Cells.Item(72,8).Value = "NORSAL"
$ws2.Cells.Item(72,9).Value = "das"
$ws2.Cells.Item(72,10).Value = "das"
$ws2.Cells.Item(72,11).Value = 68
$ws2.Cells.Item(72,12).Value = "das"
$ws2.Cells.Item(72,13).Value = "das"
$ws2.Cells.Item(72,14).Value = 265

$ws2.Cells.Item(73,1).Value = "ENTRADA"
$ws2.Cells.Item(73,2).Value = "das"
$ws2.Cells.Item(73,3).Value = "das"
$ws2.Cells.Item(73,4).Value = "das"
$ws2.Cells.Item(73,5).Value = "das"
$ws2.Cells.Item(73,6).Value = "SAL REFINADO"
$ws2.Cells.Item(73,7).Value = "Selecione uma opção"
$ws2.Cells.Item(73,8).Value = "NORSAL"
$ws2.Cells.Item(73,9).Value = "das"
$ws2.Cells.Item(73,10).Value = "das"
$ws2.Cells.Item(73,11).Value = 165
$ws2.Cells.Item(73,12).Value = "das"
$ws2.Cells.Item(73,13).Value = "das"
$ws2.Cells.Item(73,14).Value = 65

$ws2.Cells.Item(74,1).Value = "ENTRADA"
$ws2.Cells.Item(74,2).Value = "das"
$ws2.Cells.Item(74,3).Value = "das"
$ws2.Cells.Item(74,4).Value = "das"
$ws2.Cells.Item(74,5).Value = "das"
$ws2.Cells.Item(74,6).Value = "SAL REFINADO"
$ws2.Cells.Item(74,7).Value = "Selecione uma opção"
$ws2.Cells.Item(74,8).Value = "NORSAL"
$ws2.Cells.Item(74,9).Value = "das"
$ws2.Cells.Item(74,10).Value = "das"
$ws2.Cells.Item(74,11).Value = 145
$ws2.Cells.Item(74,12).Value = "das"
$ws2.Cells.Item(74,13).Value = "das"
$ws2.Cells.Item(74,14).Value = 26

# --- Sheet "Descarga do Sal": update existing cell values ---
$ws3 = $wb.Worksheets.Item("Descarga do Sal")
$ws3.Cells.Item(10,4).Value = "dasd"  # D10
$ws3.Cells.Item(20,16).Value = 356  # P20
$ws3.Cells.Item(28,4).Value = "das"  # D28
$ws3.Cells.Item(28,11).Value = 265  # K28
$ws3.Cells.Item(28,15).Value = 68  # O28
$ws3.Cells.Item(30,11).Value = 65  # K30
$ws3.Cells.Item(30,15).Value = 165  # O30
$ws3.Cells.Item(32,11).Value = 26  # K32
$ws3.Cells.Item(32,15).Value = 145  # O32

Write-Host "Done applying edits."
